$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Status column (E2:E11) from OPTIMAL to TIME_LIMIT
$ws.Range("E2:E11").Value2 = "TIME_LIMIT"

# Updated objective (B), gap (C), and solve time (D) values for fixed recourse data correction
$ws.Range("B2").Value2 = -1726.3126788491722
$ws.Range("C2").Value2 = 8.882671417141983
$ws.Range("D2").Value2 = 5599.830896778

$ws.Range("B3").Value2 = -1740.0765152806462
$ws.Range("C3").Value2 = 8.027488494657971
$ws.Range("D3").Value2 = 5523.194639918

$ws.Range("B4").Value2 = -1732.2668486303287
$ws.Range("C4").Value2 = 9.519167611773694
$ws.Range("D4").Value2 = 5488.438084554

$ws.Range("B5").Value2 = -1731.213773381394
$ws.Range("C5").Value2 = 10.14910946044649
$ws.Range("D5").Value2 = 5505.225268805

$ws.Range("B6").Value2 = -1743.143823798119
$ws.Range("C6").Value2 = 7.643987046190622
$ws.Range("D6").Value2 = 5549.293836621

$ws.Range("B7").Value2 = -1719.7767196432865
$ws.Range("C7").Value2 = 7.781098670675725
$ws.Range("D7").Value2 = 5530.258955881

$ws.Range("B8").Value2 = -1732.690555625326
$ws.Range("C8").Value2 = 8.88929688644903
$ws.Range("D8").Value2 = 5538.292544806

$ws.Range("B9").Value2 = -1716.7979290041626
$ws.Range("C9").Value2 = 9.388813050529913
$ws.Range("D9").Value2 = 5508.74755741

$ws.Range("B10").Value2 = -1730.9641896387025
$ws.Range("C10").Value2 = 7.424981119889535
$ws.Range("D10").Value2 = 5550.719144293

$ws.Range("B11").Value2 = -1743.1617425024665
$ws.Range("C11").Value2 = 8.753854342914902
$ws.Range("D11").Value2 = 5541.731702889
